# Refresh the cryptos price/volume snapshot (scheduled GitHub Actions data pull).
# Each row is "Coin (B) | Link (C) | Price (D) | Volume(1h) (E)"; column A is a
# static rank index and is left untouched. A handful of rows also changed their
# rank order, so B/C (name/link) are rewritten for those rows too.
#
# Price values that are plain decimals (e.g. "7.12") are prefixed with a leading
# apostrophe so Excel stores them as text (matching the source data's inlineStr
# cells) instead of auto-coercing them into numbers; values that already contain
# multiple dots (e.g. "56.476.54") are never auto-coerced, so no prefix is needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "56.476.54"
$ws.Range("E2").Value = "  -2.87%  "

# Row 3
$ws.Range("D3").Value = "2.943.66"

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").Value = "'494.26"
$ws.Range("E5").Value = "  -6.36%  "

# Row 6
$ws.Range("D6").Value = "'133.57"
$ws.Range("E6").Value = "  -7.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("E8").Value = "  -5.97%  "

# Row 9
$ws.Range("D9").Value = "'7.12"
$ws.Range("E9").Value = "  -6.60%  "

# Row 10
$ws.Range("E10").Value = "  -7.27%  "

# Row 11
$ws.Range("D11").Value = "'0.350"
$ws.Range("E11").Value = "  -5.74%  "

# Row 12
$ws.Range("D12").Value = "3.446.84"
$ws.Range("E12").Value = "  -3.96%  "

# Row 13
$ws.Range("E13").Value = "  -3.76%  "

# Row 14
$ws.Range("D14").Value = "'25.85"
$ws.Range("E14").Value = "  -5.67%  "

# Row 15
$ws.Range("D15").Value = "'0.0000156"
$ws.Range("E15").Value = "  -10.08%  "

# Row 16
$ws.Range("D16").Value = "56.481.80"
$ws.Range("E16").Value = "  -2.76%  "

# Row 17
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'5.96"
$ws.Range("E17").Value = "  -4.18%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.948.79"
$ws.Range("E18").Value = "  -3.60%  "

# Row 19
$ws.Range("D19").Value = "'12.42"
$ws.Range("E19").Value = "  -5.87%  "

# Row 20
$ws.Range("D20").Value = "'7.72"
$ws.Range("E20").Value = "  -6.00%  "

# Row 21
$ws.Range("D21").Value = "'316.14"
$ws.Range("E21").Value = "  -7.63%  "

# Row 22
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$ws.Range("D23").Value = "'5.75"
$ws.Range("E23").Value = "  +0.26%  "

# Row 24
$ws.Range("D24").Value = "'0.483"
$ws.Range("E24").Value = "  -4.50%  "

# Row 25
$ws.Range("D25").Value = "'62.31"
$ws.Range("E25").Value = "  -4.88%  "

# Row 26
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.33%  "

# Row 27
$ws.Range("E27").Value = "  -5.67%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0852"
$ws.Range("E28").Value = "  -12.74%  "

# Row 29
$ws.Range("D29").Value = "'6.42"
$ws.Range("E29").Value = "  -8.23%  "

# Row 30
$ws.Range("D30").Value = "'7.00"
$ws.Range("E30").Value = "  -6.63%  "

# Row 31
$ws.Range("E31").Value = "  -7.02%  "

# Row 32
$ws.Range("D32").Value = "'19.91"
$ws.Range("E32").Value = "  -6.15%  "

# Row 33
$ws.Range("E33").Value = "  -8.90%  "

# Row 34
$ws.Range("D34").Value = "'151.51"
$ws.Range("E34").Value = "  -3.95%  "

# Row 35
$ws.Range("D35").Value = "'4.44"
$ws.Range("E35").Value = "  -8.08%  "

# Row 36
$ws.Range("D36").Value = "'5.65"
$ws.Range("E36").Value = "  -5.65%  "

# Row 37
$ws.Range("D37").Value = "'1.20"
$ws.Range("E37").Value = "  -9.96%  "

# Row 38
$ws.Range("D38").Value = "'23.59"
$ws.Range("E38").Value = "  -10.15%  "

# Row 39
$ws.Range("D39").Value = "'0.0651"
$ws.Range("E39").Value = "  -7.07%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'37.31"
$ws.Range("E40").Value = "  -1.45%  "

# Row 41
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "2.966.02"
$ws.Range("E41").Value = "  -4.28%  "

# Row 42
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  -0.22%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.66"
$ws.Range("E43").Value = "  -8.36%  "

# Row 44
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.632"
$ws.Range("E44").Value = "  -5.59%  "

# Row 45
$ws.Range("D45").Value = "2.133.07"
$ws.Range("E45").Value = "  -8.53%  "

# Row 46
$ws.Range("D46").Value = "'1.33"
$ws.Range("E46").Value = "  -10.22%  "

# Row 47
$ws.Range("D47").Value = "'5.84"
$ws.Range("E47").Value = "  -4.23%  "

# Row 48
$ws.Range("D48").Value = "'0.909"
$ws.Range("E48").Value = "  -13.09%  "

# Row 49
$ws.Range("D49").Value = "'0.0229"
$ws.Range("E49").Value = "  -6.39%  "

# Row 50
$ws.Range("D50").Value = "'18.86"
$ws.Range("E50").Value = "  -7.20%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0847"
$ws.Range("E51").Value = "  -6.55%  "
